$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cell D1, matching the style of the existing header row (A1:C1)
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("D1").Value = "Tipo"

# Update the existing numeric values in row 2
$ws.Range("B2").Value = 0.4125991217944497
$ws.Range("C2").Value = 0.9918641102167772

# Add the new value in D2
$ws.Range("D2").Value = "single"
